# Adding search test cases
# Adds two new rows (28, 29) to the "Test Cases" sheet for the new
# ProfileTypeaheadCountySelectTest / ProfileCountryTypeaheadOptionsDisplayTest
# search test cases, mirrors formatting from the last existing data row (27),
# and updates the Results column so the new final row carries the PASS
# status previously held by row 27 (which becomes SKIP).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Row 28: ProfileTypeaheadCountySelectTest -----------------------------

$ws.Range("A27").Copy()
$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("B27").Copy()
$ws.Range("B28").PasteSpecial(-4122)
$ws.Range("C27").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D27").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E27").Copy()
$ws.Range("E28").PasteSpecial(-4122)

$ws.Range("A28").Value = "ProfileTypeaheadCountySelectTest"
$ws.Range("B28").Value = "TBD"
$ws.Range("C28").Value = "Verify that user is able to add 'country' using typeahead"
$ws.Range("D28").Value = "Y"
$ws.Range("E28").Value = "SKIP"

# Italicise "country'" within the description (matches the author's rich text)
$run2 = $ws.Range("C28").Characters(34, 8)
$run2.Font.Italic = $true
$run3 = $ws.Range("C28").Characters(42, 16)
$run3.Font.Bold = $false

# --- Row 29: ProfileCountryTypeaheadOptionsDisplayTest ---------------------

$ws.Range("A27").Copy()
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("B27").Copy()
$ws.Range("B29").PasteSpecial(-4122)
$ws.Range("C27").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("D27").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E27").Copy()
$ws.Range("E29").PasteSpecial(-4122)

$ws.Range("A29").Value = "ProfileCountryTypeaheadOptionsDisplayTest"
$ws.Range("B29").Value = "TBD"
$ws.Range("C29").Value = "Verify that  'country' using type ahead options should display while enter min 2 characters"
$ws.Range("D29").Value = "Y"
$ws.Range("E29").Value = "PASS"

# Italicise the single "'" character within the description
$run2b = $ws.Range("C29").Characters(22, 1)
$run2b.Font.Italic = $true
$run3b = $ws.Range("C29").Characters(23, 69)
$run3b.Font.Bold = $false

# Row 27 no longer holds the last PASS result - it becomes SKIP now that
# row 29 is the new final (PASS) row.
$ws.Range("E27").Value = "SKIP"

# --- Sheet view bookkeeping --------------------------------------------

$ws.Range("D29").Select()
